$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original sheet had an empty column E (between "notes" and "SiO2") and an
# empty column Q (between "K2O" and "tot"). Remove both columns, which shifts
# the remaining data left so SiO2..tot now occupy E..P.
$ws.Columns("E").Delete() | Out-Null
$ws.Columns("P").Delete() | Out-Null

# Leave the selection on the (now) last column, matching the final view state.
$ws.Columns("P").Select() | Out-Null
